# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Wed Jan  3 04:43:35 UTC 2024 with GitHub Actions"
#
# The sheet lists crypto coins with columns: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h).
# All cells are stored as text. Column D sometimes holds values that look like plain
# numbers (e.g. "8.47"), so we force the cell's number format to Text ("@") before
# writing so Excel does not silently reinterpret the string as a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.328.49'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.367.65'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.53'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.92'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.76'
$ws.Range("E10").Value = '  -3.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0916'
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.47'
$ws.Range("E12").Value = '  -2.37%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  -3.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.727.90'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.32'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.356.98'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.296.69'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.41'
$ws.Range("E19").Value = '  +17.68%  '
$ws.Range("E20").Value = '  -5.29%  '
$ws.Range("E21").Value = '  -2.11%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.29'
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.09'
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.49'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.15'
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0968'
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.29'
$ws.Range("E31").Value = '  -2.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.02'
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.54'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.88'
$ws.Range("E34").Value = '  -4.02%  '
$ws.Range("E35").Value = '  -2.13%  '
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.69'
$ws.Range("E37").Value = '  -3.97%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.97'
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("E39").Value = '  +7.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("E40").Value = '  -5.52%  '
$ws.Range("E41").Value = '  -2.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.57'
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.13'
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.04'
$ws.Range("E44").Value = '  -1.65%  '
$ws.Range("E45").Value = '  -5.60%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.820.07'
$ws.Range("E48").Value = '  +10.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '82.63'
$ws.Range("E49").Value = '  +4.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.06'
$ws.Range("E50").Value = '  -5.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.18'
$ws.Range("E51").Value = '  -1.85%  '
